$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.4960329532623291
$ws.Range("B1").Value = 1.33648145198822
$ws.Range("C1").Value = 4.870220184326172
$ws.Range("D1").Value = 1.434517979621887
$ws.Range("E1").Value = 0.8113567233085632
